# Rename "Fixation_*" filenames to "Ball_*" filenames in the Filename_Left /
# Filename_Right columns (D and E). These cells were mislabeled as fixation
# images when they actually refer to the ball stimulus images.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rename = @{
    "Fixation_D64_L.png" = "Ball_D64_L.png"
    "Fixation_D64_R.png" = "Ball_D64_R.png"
    "Fixation_D51_L.png" = "Ball_D51_L.png"
    "Fixation_D51_R.png" = "Ball_D51_R.png"
    "Fixation_D80_L.png" = "Ball_D80_L.png"
    "Fixation_D80_R.png" = "Ball_D80_R.png"
}

$lastRow = 147

# Process column E (Filename_Right) fully top-to-bottom first, then column D
# (Filename_Left), so that the new shared-string entries get appended in the
# same order as the authored edit (R-variants before L-variants).
foreach ($col in 5, 4) {
    for ($r = 2; $r -le $lastRow; $r++) {
        $cell = $ws.Cells.Item($r, $col)
        $val = $cell.Value2
        if ($rename.ContainsKey($val)) {
            $cell.Value = $rename[$val]
        }
    }
}

# Update the selection to match the authored change.
$ws.Range("F5").Select()
